$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.57
$ws.Range("G2").Value = 1.68
$ws.Range("H2").Value = 5.4
$ws.Range("J2").Value = 4.1
$ws.Range("P2").Value = 2.12
$ws.Range("T2").Value = 1.73
$ws.Range("U2").Value = 2.1
$ws.Range("X2").Value = 22
$ws.Range("Y2").Value = 26
$ws.Range("Z2").Value = 60
$ws.Range("AA2").Value = 190
$ws.Range("AB2").Value = 11
$ws.Range("AC2").Value = 11
$ws.Range("AD2").Value = 25
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 11.5
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 20
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 16.5
$ws.Range("AK2").Value = 16.5
$ws.Range("AL2").Value = 32
$ws.Range("AM2").Value = 110
$ws.Range("AN2").Value = 7.6
$ws.Range("AO2").Value = 95

# Row 3
$ws.Range("F3").Value = 2.28
$ws.Range("I3").Value = 3.5
$ws.Range("J3").Value = 3.5
$ws.Range("M3").Value = 1.05
$ws.Range("Q3").Value = 1.73
$ws.Range("U3").Value = 2.28
$ws.Range("Z3").Value = 26
$ws.Range("AD3").Value = 17
$ws.Range("AE3").Value = 42
$ws.Range("AF3").Value = 17.5
$ws.Range("AH3").Value = 19.5
$ws.Range("AI3").Value = 50
$ws.Range("AJ3").Value = 32
$ws.Range("AK3").Value = 25
$ws.Range("AM3").Value = 90
$ws.Range("AN3").Value = 19.5

# Row 4
$ws.Range("I4").Value = 27
$ws.Range("J4").Value = 6.4
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 2.78
$ws.Range("O4").Value = 1.13
$ws.Range("Q4").Value = 1.55
$ws.Range("R4").Value = 1.32
$ws.Range("S4").Value = 2.34
$ws.Range("T4").Value = 1.01
$ws.Range("U4").Value = 1.01
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 1000

# Row 5
$ws.Range("F5").Value = 1.64
$ws.Range("G5").Value = 1.71
$ws.Range("H5").Value = 5.9
$ws.Range("J5").Value = 3.9
$ws.Range("K5").Value = 4.5
$ws.Range("M5").Value = 1.07
$ws.Range("P5").Value = 1.75
$ws.Range("Q5").Value = 1.91
$ws.Range("T5").Value = 2.1
$ws.Range("U5").Value = 1.73

# Row 6
$ws.Range("N6").Value = 3.55

# Row 7
$ws.Range("F7").Value = 1.3
$ws.Range("G7").Value = 1.36
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 13.5
$ws.Range("AF7").Value = 8.6
$ws.Range("AM7").Value = 210

# Row 8
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 4.6
$ws.Range("J8").Value = 3.4
$ws.Range("K8").Value = 3.8

# Row 9
$ws.Range("G9").Value = 2.82
$ws.Range("H9").Value = 2.76
$ws.Range("J9").Value = 3.3
$ws.Range("P9").Value = 1.68
$ws.Range("Q9").Value = 1.86
$ws.Range("U9").Value = 2
$ws.Range("AB9").Value = 12.5
$ws.Range("AC9").Value = 9
$ws.Range("AD9").Value = 16
$ws.Range("AF9").Value = 21
$ws.Range("AG9").Value = 14.5

# Row 10
$ws.Range("F10").Value = 3.45
$ws.Range("G10").Value = 3.5
$ws.Range("H10").Value = 2.78
$ws.Range("I10").Value = 2.82
$ws.Range("J10").Value = 2.8
$ws.Range("K10").Value = 2.82
$ws.Range("P10").Value = 1.26
$ws.Range("Q10").Value = 1.01
$ws.Range("X10").Value = 1000

# Row 11
$ws.Range("G11").Value = 1.3
$ws.Range("H11").Value = 8.4
$ws.Range("I11").Value = 14
$ws.Range("P11").Value = 2.72
$ws.Range("Q11").Value = 1.39
$ws.Range("T11").Value = 1.88
$ws.Range("AB11").Value = 16
$ws.Range("AC11").Value = 21
$ws.Range("AG11").Value = 13
$ws.Range("AJ11").Value = 13
$ws.Range("AK11").Value = 15
$ws.Range("AN11").Value = 4.1

# Row 12
$ws.Range("F12").Value = 1.31
$ws.Range("G12").Value = 1.37
$ws.Range("H12").Value = 9
$ws.Range("J12").Value = 5.6
$ws.Range("K12").Value = 6.6
$ws.Range("P12").Value = 2.34
$ws.Range("Q12").Value = 1.58
$ws.Range("T12").Value = 2.02
$ws.Range("U12").Value = 1.79
$ws.Range("AE12").Value = 220
$ws.Range("AI12").Value = 170
$ws.Range("AK12").Value = 17.5
$ws.Range("AN12").Value = 5.8

# Row 13
$ws.Range("F13").Value = 1.47
$ws.Range("G13").Value = 1.5
$ws.Range("H13").Value = 7.2
$ws.Range("I13").Value = 11
$ws.Range("J13").Value = 4.4
$ws.Range("K13").Value = 5.3
$ws.Range("P13").Value = 1.9
$ws.Range("Q13").Value = 1.87
$ws.Range("T13").Value = 2.06
$ws.Range("U13").Value = 1.71
$ws.Range("X13").Value = 19
$ws.Range("AA13").Value = 500
$ws.Range("AL13").Value = 55

# Row 14
$ws.Range("H14").Value = 15
$ws.Range("J14").Value = 7.6
$ws.Range("T14").Value = 2.36
$ws.Range("Z14").Value = 360
$ws.Range("AJ14").Value = 8.800000000000001
$ws.Range("AK14").Value = 16
$ws.Range("AM14").Value = 340
$ws.Range("AN14").Value = 3.6

# Row 15
$ws.Range("F15").Value = 2.12
$ws.Range("I15").Value = 3.9
$ws.Range("J15").Value = 3.7

# Row 16
$ws.Range("H16").Value = 10
$ws.Range("K16").Value = 7.2
$ws.Range("P16").Value = 2.88
$ws.Range("Q16").Value = 1.4
$ws.Range("T16").Value = 1.89
$ws.Range("U16").Value = 1.92
$ws.Range("X16").Value = 44
$ws.Range("Y16").Value = 1000
$ws.Range("Z16").Value = 150
$ws.Range("AA16").Value = 500
$ws.Range("AB16").Value = 15
$ws.Range("AF16").Value = 12
$ws.Range("AH16").Value = 1000
$ws.Range("AM16").Value = 150
$ws.Range("AN16").Value = 4.2

# Row 17
$ws.Range("G17").Value = 3.55
$ws.Range("H17").Value = 2.22
$ws.Range("I17").Value = 2.24
$ws.Range("J17").Value = 3.65
$ws.Range("P17").Value = 2
$ws.Range("Q17").Value = 1.78
$ws.Range("T17").Value = 1.57
$ws.Range("U17").Value = 2.22
$ws.Range("Y17").Value = 11.5
$ws.Range("AD17").Value = 11.5
$ws.Range("AK17").Value = 44
$ws.Range("AO17").Value = 16

# Row 18
$ws.Range("F18").Value = 1.91
$ws.Range("G18").Value = 2.04
$ws.Range("H18").Value = 4
$ws.Range("I18").Value = 4.7
$ws.Range("J18").Value = 3.5
$ws.Range("M18").Value = 1.07
$ws.Range("P18").Value = 1.8
$ws.Range("Q18").Value = 1.85
$ws.Range("T18").Value = 1.85
$ws.Range("U18").Value = 1.96
$ws.Range("X18").Value = 16.5
$ws.Range("AB18").Value = 1000
$ws.Range("AC18").Value = 1000
$ws.Range("AM18").Value = 150

# Row 19
$ws.Range("F19").Value = 1.71
$ws.Range("G19").Value = 1.77
$ws.Range("I19").Value = 6.4
$ws.Range("K19").Value = 4.5
$ws.Range("P19").Value = 2.08
$ws.Range("AJ19").Value = 980
